# Adds two new slides ("Jodel" and "Jodel Alert") describing the domain,
# using the same "Title and Content" layout as used elsewhere in the deck.

$p = $ppt.ActivePresentation

function Set-JodelFont($range) {
    $range.Font.Name = "Arial Unicode MS"
    $range.Font.NameFarEast = "Arial Unicode MS"
    $range.Font.NameComplexScript = "Arial Unicode MS"
}

# ---------------------------------------------------------------------
# Slide 3: "Jodel"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Add(3, 2)

$s3.Shapes.Item(1).Name = "Title 1"
$s3.Shapes.Item(2).Name = "Content Placeholder 2"

$title3 = $s3.Shapes.Item(1).TextFrame.TextRange
$title3.Text = "Jodel"
$title3.LanguageID = "sv-SE"
Set-JodelFont $title3

$body3 = $s3.Shapes.Item(2).TextFrame.TextRange
$body3.Text = "An anonymous social application that targets students and campus life"
$body3.LanguageID = "sv-SE"
Set-JodelFont $body3

$body3Empty = $body3.InsertAfter("`r`r")
$body3Empty.LanguageID = "sv-SE"
Set-JodelFont $body3Empty

# ---------------------------------------------------------------------
# Slide 4: "Jodel Alert"
# ---------------------------------------------------------------------
$s4 = $p.Slides.Add(4, 2)

$s4.Shapes.Item(1).Name = "Title 1"
$s4.Shapes.Item(2).Name = "Content Placeholder 2"

$title4 = $s4.Shapes.Item(1).TextFrame.TextRange
$title4.Text = "Jodel Alert"
$title4.LanguageID = "sv-SE"
Set-JodelFont $title4

$body4 = $s4.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1 (two runs)
$body4.Text = "Linnéstudenterna wants to tap into the feed and get an alert when someone posts about accomodation, exams, cheating and students getting mistreated from the university "
$body4.LanguageID = "sv-SE"
Set-JodelFont $body4
$p1r2 = $body4.InsertAfter("etc (keywords)")
$p1r2.LanguageID = "sv-SE"
Set-JodelFont $p1r2

# Paragraph 2 (empty)
$p2 = $p1r2.InsertAfter("`r`r")
$p2.LanguageID = "sv-SE"
Set-JodelFont $p2

# Paragraph 3 (three runs)
$p3r1 = $p2.InsertAfter("Jodel Alert will ")
$p3r1.LanguageID = "sv-SE"
Set-JodelFont $p3r1
$p3r2 = $p3r1.InsertAfter("send Linnéstudenterna an email when such post is ")
$p3r2.LanguageID = "sv-SE"
Set-JodelFont $p3r2
$p3r3 = $p3r2.InsertAfter("found")
$p3r3.LanguageID = "sv-SE"
Set-JodelFont $p3r3

# Paragraph 4 (empty)
$p4 = $p3r3.InsertAfter("`r`r")
$p4.LanguageID = "sv-SE"
Set-JodelFont $p4

# Paragraph 5
$p5 = $p4.InsertAfter("Today Linnéstudenterna searches the feed manually when given time")
$p5.LanguageID = "sv-SE"
Set-JodelFont $p5

# Paragraph 6 (empty)
$p6 = $p5.InsertAfter("`r`r")
$p6.LanguageID = "sv-SE"
Set-JodelFont $p6

# Paragraph 7
$p7 = $p6.InsertAfter("Other companies can also benefit from this, in order to see what is trending for that company or what is said in general public by common man about that company")
$p7.LanguageID = "sv-SE"
Set-JodelFont $p7
